# mals2-44 - address + formatting changes
#
# 1) "Minister of Finance" gets its own (slightly smaller, 10.5pt) run so it
#    can be emphasised independently of the surrounding bold/italic/red text.
# 2) "Livestock Health Management and Regulation" (the sign-off block near
#    the mailing address) is renamed to "Office of the Chief Veterinarian".
# 3) The mailing address is tidied up ("B.C." -> "BC", drop the trailing
#    space/tab) and the now-stale Telephone / Toll-Free lines (and the
#    blank paragraphs around them) are removed from the footer entirely.

$d = $word.ActiveDocument

# --- 1) Split out "Minister of Finance" into its own run at 10.5pt -------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Minister of Finance", $false, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # rng now covers just "Minister of Finance"; bumping the size splits
    # it away from the " " runs on either side, which keep their original
    # (unset / inherited) size.
    $rng.Font.Size = 10.5
}

# --- 2) Rename the division line under "Ministry of Agriculture and Food" -
# (There are two "Livestock Health Management and Regulation" paragraphs in
# this template - a centered heading near the top of the page, and this
# bold footer line just above the mailing address. Only the footer one
# changes, so match on its distinguishing "both"-justified style by
# looking at the paragraph immediately preceding the address block.)
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "Livestock Health Management and Regulation" -and `
        $p.Range.ParagraphFormat.Alignment -eq 3) {
        $pr = $p.Range
        $pr.MoveEnd(1, -1) | Out-Null
        $pr.Text = "Office of the Chief Veterinarian"
        break
    }
}

# --- 3) Clean up the address block / drop the phone-number paragraphs ----
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Abbotsford, B.C.*") {
        $addrIndex = $i
        $ar = $p.Range
        $ar.MoveEnd(1, -1) | Out-Null
        $ar.Text = "Abbotsford, BC   V3G 2M3"
        break
    }
}

if ($addrIndex) {
    # The paragraph right after the address is a blank line, then
    # "Telephone: ...", then "Toll-Free (BC) ...", then a trailing blank
    # paragraph. All four are being removed; deleting the combined range
    # collapses them away and leaves the address paragraph followed
    # directly by the section break.
    $startPara = $d.Paragraphs.Item($addrIndex + 1)
    $endPara = $d.Paragraphs.Item($addrIndex + 4)
    $killRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $killRange.Delete()
}
